# Insert a new data row before the current row 569 (shifts existing rows 569-607 down to 570-608)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(569).Insert()

$ws.Range("A569").Value = 10
$ws.Range("B569").Value = "Vega Modelo de Temuco"
$ws.Range("C569").Value = "La Araucanía"
$ws.Range("D569").Value = 44826
$ws.Range("E569").Value = 9
$ws.Range("F569").Value = 100112027
$ws.Range("G569").Value = "Melón"
$ws.Range("H569").Value = "Calameño"
$ws.Range("I569").Value = "Extra"
$ws.Range("J569").Value = 80
$ws.Range("K569").Value = 27000
$ws.Range("L569").Value = 27000
$ws.Range("M569").Value = 27000
$ws.Range("N569").Value = "$/caja 12 unidades"
$ws.Range("O569").Value = "Brasil"
$ws.Range("P569").Value = 2250
$ws.Range("Q569").Value = 12
$ws.Range("R569").Value = "Hortaliza"
